# "Generate Report for Handoff"
#
# Two new localizable source files showed up since the last report:
#   9015650a-e053-4466-aeaa-9d9227d4ceb1.md
#   b2f489bc-9159-4fc1-9de3-126ab88b431f.md
# Both are "Ready for handoff" and already have a first handoff record
# (zh-cn + de-de). The sentinel ".localization-config" row (not localized)
# simply shifts down to make room for them.

$wb = $excel.ActiveWorkbook

$mdBase  = "https://github.com/OpenLocalizationTest/oltest/blob/89b067407e92bddf2e3ac65baad03bfe79345652"
$zhBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e10d945f27b34c2b08533ff306d7369e460faaec/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht"
$deBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33a7094fd424c430033739bdfde01a7bd58219f3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht"

$file1     = "9015650a-e053-4466-aeaa-9d9227d4ceb1.md"
$file1zh   = "9015650a-e053-4466-aeaa-9d9227d4ceb1.00e409ac0edf6896bc2ad19456cebb33cde08449.zh-cn.xlf"
$file1de   = "9015650a-e053-4466-aeaa-9d9227d4ceb1.00e409ac0edf6896bc2ad19456cebb33cde08449.de-de.xlf"
$file1time = "2016-03-09 07:46:32"
$file1timeDe = "2016-03-09 07:46:44"

$file2     = "b2f489bc-9159-4fc1-9de3-126ab88b431f.md"
$file2zh   = "b2f489bc-9159-4fc1-9de3-126ab88b431f.dc341845429dddf0be9c7d5c99e5a83e3f3c80e9.zh-cn.xlf"
$file2de   = "b2f489bc-9159-4fc1-9de3-126ab88b431f.dc341845429dddf0be9c7d5c99e5a83e3f3c80e9.de-de.xlf"

$readyForHandoff = "Ready for handoff"
$notLocalized    = "Not to be localized"
$cfgName         = ".localization-config"
$epoch           = "0001-01-01 00:00:00"

# ---------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("A4").Value = $file1
$ws1.Range("B4").Value = $readyForHandoff
$ws1.Range("C4").Value = $readyForHandoff

$ws1.Range("A5").Value = $file2
$ws1.Range("B5").Value = $readyForHandoff
$ws1.Range("C5").Value = $readyForHandoff

$ws1.Range("A6").Value = $cfgName
$ws1.Range("B6").Value = $notLocalized
$ws1.Range("C6").Value = $notLocalized

$ws1.Hyperlinks.Add($ws1.Range("A2"), "$mdBase/e2e/0ab8bbd2-1d59-44ca-8cdf-bff9e35b7ae2.md", "", "", "0ab8bbd2-1d59-44ca-8cdf-bff9e35b7ae2.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$mdBase/e2e/2a9e401f-3f63-48d3-9a02-5965d6b8c6b9.md", "", "", "2a9e401f-3f63-48d3-9a02-5965d6b8c6b9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$mdBase/e2e/$file1", "", "", $file1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "$mdBase/e2e/$file2", "", "", $file2) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "$mdBase/$cfgName", "", "", $cfgName) | Out-Null

# ---------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("A4").Value = $file1
$ws2.Range("B4").Value = $readyForHandoff
$ws2.Range("C4").Value = $file1zh
$ws2.Range("D4").Value = $file1time
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = "Include"

$ws2.Range("A5").Value = $file2
$ws2.Range("B5").Value = $readyForHandoff
$ws2.Range("C5").Value = $file2zh
$ws2.Range("D5").Value = $file1time
$ws2.Range("G5").Value = $epoch
$ws2.Range("H5").Value = "Include"

$ws2.Range("A6").Value = $cfgName
$ws2.Range("B6").Value = $notLocalized
$ws2.Range("D6").Value = $epoch
$ws2.Range("G6").Value = $epoch
$ws2.Range("H6").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "$mdBase/e2e/0ab8bbd2-1d59-44ca-8cdf-bff9e35b7ae2.md", "", "", "0ab8bbd2-1d59-44ca-8cdf-bff9e35b7ae2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhBase/0ab8bbd2-1d59-44ca-8cdf-bff9e35b7ae2.17616d05980ad4cf730d8c598a988901bc846889.zh-cn.xlf", "", "", "0ab8bbd2-1d59-44ca-8cdf-bff9e35b7ae2.17616d05980ad4cf730d8c598a988901bc846889.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$mdBase/e2e/2a9e401f-3f63-48d3-9a02-5965d6b8c6b9.md", "", "", "2a9e401f-3f63-48d3-9a02-5965d6b8c6b9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$zhBase/2a9e401f-3f63-48d3-9a02-5965d6b8c6b9.fd784d34a0e0b5f9b18cea3f002bddd7cc8a5846.zh-cn.xlf", "", "", "2a9e401f-3f63-48d3-9a02-5965d6b8c6b9.fd784d34a0e0b5f9b18cea3f002bddd7cc8a5846.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$mdBase/e2e/$file1", "", "", $file1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "$zhBase/$file1zh", "", "", $file1zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "$mdBase/e2e/$file2", "", "", $file2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "$zhBase/$file2zh", "", "", $file2zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "$mdBase/$cfgName", "", "", $cfgName) | Out-Null

# ---------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("A4").Value = $file1
$ws3.Range("B4").Value = $readyForHandoff
$ws3.Range("C4").Value = $file1de
$ws3.Range("D4").Value = $file1timeDe
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = "Include"

$ws3.Range("A5").Value = $file2
$ws3.Range("B5").Value = $readyForHandoff
$ws3.Range("C5").Value = $file2de
$ws3.Range("D5").Value = $file1timeDe
$ws3.Range("G5").Value = $epoch
$ws3.Range("H5").Value = "Include"

$ws3.Range("A6").Value = $cfgName
$ws3.Range("B6").Value = $notLocalized
$ws3.Range("D6").Value = $epoch
$ws3.Range("G6").Value = $epoch
$ws3.Range("H6").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "$mdBase/e2e/0ab8bbd2-1d59-44ca-8cdf-bff9e35b7ae2.md", "", "", "0ab8bbd2-1d59-44ca-8cdf-bff9e35b7ae2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deBase/0ab8bbd2-1d59-44ca-8cdf-bff9e35b7ae2.17616d05980ad4cf730d8c598a988901bc846889.de-de.xlf", "", "", "0ab8bbd2-1d59-44ca-8cdf-bff9e35b7ae2.17616d05980ad4cf730d8c598a988901bc846889.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$mdBase/e2e/2a9e401f-3f63-48d3-9a02-5965d6b8c6b9.md", "", "", "2a9e401f-3f63-48d3-9a02-5965d6b8c6b9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$deBase/2a9e401f-3f63-48d3-9a02-5965d6b8c6b9.fd784d34a0e0b5f9b18cea3f002bddd7cc8a5846.de-de.xlf", "", "", "2a9e401f-3f63-48d3-9a02-5965d6b8c6b9.fd784d34a0e0b5f9b18cea3f002bddd7cc8a5846.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$mdBase/e2e/$file1", "", "", $file1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "$deBase/$file1de", "", "", $file1de) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "$mdBase/e2e/$file2", "", "", $file2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "$deBase/$file2de", "", "", $file2de) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "$mdBase/$cfgName", "", "", $cfgName) | Out-Null
